$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell values (recalculated stats) ---
$ws.Range("G9").Value = 0.27923518304623
$ws.Range("G10").Value = 0.27923518304623
$ws.Range("G11").Value = 0.869977739057868
$ws.Range("G12").Value = 0.869977739057868
$ws.Range("F13").Value = 0.4042
$ws.Range("G13").Value = 1.33524210526316
$ws.Range("F14").Value = 0.4042
$ws.Range("G14").Value = 1.33524210526316
$ws.Range("G26").Value = 0.20249949664586
$ws.Range("G27").Value = 0.20249949664586
$ws.Range("G28").Value = 0.72562024097754
$ws.Range("G29").Value = 0.72562024097754
$ws.Range("F30").Value = 0.2029
$ws.Range("G30").Value = 1.07615614035088
$ws.Range("L30").Value = 0.0505
$ws.Range("F31").Value = 0.2029
$ws.Range("G31").Value = 1.07615614035088
$ws.Range("L31").Value = 0.0505
$ws.Range("G36").Value = 1.46228762606776
$ws.Range("G43").Value = 0.186765977223301
$ws.Range("G44").Value = 0.186765977223301
$ws.Range("G45").Value = 0.40996519537765
$ws.Range("L45").Value = 0.00685
$ws.Range("G46").Value = 0.40996519537765
$ws.Range("L46").Value = 0.00685
$ws.Range("F47").Value = 0.2151
$ws.Range("G47").Value = 0.72174767834455
$ws.Range("F48").Value = 0.2151
$ws.Range("G48").Value = 0.72174767834455
$ws.Range("G53").Value = 1.44436735841827
$ws.Range("G60").Value = 0.0963153010427784
$ws.Range("L60").Value = 0.00324
$ws.Range("G61").Value = 0.0963153010427784
$ws.Range("L61").Value = 0.00324
$ws.Range("G62").Value = 0.17823758582355
$ws.Range("L62").Value = 0.00351
$ws.Range("G63").Value = 0.17823758582355
$ws.Range("L63").Value = 0.00351
$ws.Range("F64").Value = 0.1471
$ws.Range("G64").Value = 0.343223679762748
$ws.Range("L64").Value = 0.01665
$ws.Range("M64").Value = 0.62418
$ws.Range("F65").Value = 0.1471
$ws.Range("G65").Value = 0.343223679762748
$ws.Range("L65").Value = 0.01665
$ws.Range("M65").Value = 0.62418
$ws.Range("G70").Value = 1.26976165888516
$ws.Range("F77").Value = 0.00472
$ws.Range("G77").Value = 0.0500634001498791
$ws.Range("L77").Value = 0.00276
$ws.Range("F78").Value = 0.00472
$ws.Range("G78").Value = 0.0500634001498791
$ws.Range("L78").Value = 0.00276
$ws.Range("G79").Value = 0.13670220552921
$ws.Range("L79").Value = 0.00565
$ws.Range("G80").Value = 0.13670220552921
$ws.Range("L80").Value = 0.00565
$ws.Range("F81").Value = 0.16215
$ws.Range("G81").Value = 0.247385346429415
$ws.Range("I81").Value = 0.77815
$ws.Range("L81").Value = 0.01415
$ws.Range("M81").Value = 0.45938
$ws.Range("F82").Value = 0.16215
$ws.Range("G82").Value = 0.247385346429415
$ws.Range("I82").Value = 0.77815
$ws.Range("L82").Value = 0.01415
$ws.Range("M82").Value = 0.45938
$ws.Range("G87").Value = 1.04713901737573
$ws.Range("G94").Value = 0.0492994202169429
$ws.Range("L94").Value = 0.00482
$ws.Range("G95").Value = 0.0492994202169429
$ws.Range("L95").Value = 0.00482
$ws.Range("G96").Value = 0.14692220552921
$ws.Range("G97").Value = 0.14692220552921
$ws.Range("F98").Value = 0.17315
$ws.Range("G98").Value = 0.245657013096082
$ws.Range("F99").Value = 0.17315
$ws.Range("G99").Value = 0.245657013096082
$ws.Range("G104").Value = 0.972879416743886
$ws.Range("F111").Value = 0.01317
$ws.Range("G111").Value = 0.107228310995635
$ws.Range("L111").Value = 0.00856
$ws.Range("F112").Value = 0.01317
$ws.Range("G112").Value = 0.107228310995635
$ws.Range("L112").Value = 0.00856
$ws.Range("G113").Value = 0.253915538862544
$ws.Range("G114").Value = 0.253915538862544
$ws.Range("G115").Value = 0.383422013096082
$ws.Range("G116").Value = 0.383422013096082
$ws.Range("G128").Value = 0.110830792478832
$ws.Range("L128").Value = 0.00967
$ws.Range("G129").Value = 0.110830792478832
$ws.Range("L129").Value = 0.00967
$ws.Range("G130").Value = 0.319564552917485
$ws.Range("G131").Value = 0.319564552917485
$ws.Range("G132").Value = 0.44999674175316
$ws.Range("G133").Value = 0.44999674175316
$ws.Range("G145").Value = 0.197869001188005
$ws.Range("L145").Value = 0.01274
$ws.Range("G146").Value = 0.197869001188005
$ws.Range("L146").Value = 0.01274
$ws.Range("G149").Value = 0.586228813559322
$ws.Range("G150").Value = 0.586228813559322
$ws.Range("G162").Value = 0.2119103107103
$ws.Range("G163").Value = 0.2119103107103
$ws.Range("G166").Value = 0.610318965517241
$ws.Range("L166").Value = 0.5405
$ws.Range("G167").Value = 0.610318965517241
$ws.Range("L167").Value = 0.5405

# --- Append new rows 172-188 ---
# Row 172
$ws.Range("A172").Value = "Rangitikei at d/s Riverlands"
$ws.Range("B172").Value = "Visual Clarity (Sediment class 1)"
$ws.Range("C172").Value = "D"
$ws.Range("D172").Value = "2019 - 2023"
$ws.Range("E172").Value = "Impact"
$ws.Range("F172").Value = 0.7
$ws.Range("G172").Value = 1.11330769230769
$ws.Range("H172").Value = 3.7
$ws.Range("I172").Value = 3.27
$ws.Range("J172").Value = ""
$ws.Range("K172").Value = ""
$ws.Range("L172").Value = 1.8
$ws.Range("M172").Value = 2.066
$ws.Range("N172").Value = 2.802
$ws.Range("O172").Value = 1799980.004
$ws.Range("P172").Value = 5547896.885
$ws.Range("Q172").Value = "Rangitikei District"
$ws.Range("R172").Value = "Rangitīkei-Turakina"
$ws.Range("S172").Value = "Coastal Rangitikei"
$ws.Range("T172").Value = "Rang_4d"
$ws.Range("U172").Value = "m"

# Row 173
$ws.Range("A173").Value = "Rangitikei at d/s Riverlands"
$ws.Range("B173").Value = "DRP (95th Percentile)"
$ws.Range("C173").Value = "D"
$ws.Range("D173").Value = "2019 - 2023"
$ws.Range("E173").Value = "Impact"
$ws.Range("F173").Value = 0.07049999999999999
$ws.Range("G173").Value = 0.145948275862069
$ws.Range("H173").Value = 1.05
$ws.Range("I173").Value = 0.4666
$ws.Range("J173").Value = ""
$ws.Range("K173").Value = ""
$ws.Range("L173").Value = 0.136
$ws.Range("M173").Value = 0.30156
$ws.Range("N173").Value = 0.42132
$ws.Range("O173").Value = 1799980.004
$ws.Range("P173").Value = 5547896.885
$ws.Range("Q173").Value = "Rangitikei District"
$ws.Range("R173").Value = "Rangitīkei-Turakina"
$ws.Range("S173").Value = "Coastal Rangitikei"
$ws.Range("T173").Value = "Rang_4d"
$ws.Range("U173").Value = "mg/L"

# Row 174
$ws.Range("A174").Value = "Rangitikei at d/s Riverlands"
$ws.Range("B174").Value = "DRP (Median)"
$ws.Range("C174").Value = "D"
$ws.Range("D174").Value = "2019 - 2023"
$ws.Range("E174").Value = "Impact"
$ws.Range("F174").Value = 0.07049999999999999
$ws.Range("G174").Value = 0.145948275862069
$ws.Range("H174").Value = 1.05
$ws.Range("I174").Value = 0.4666
$ws.Range("J174").Value = ""
$ws.Range("K174").Value = ""
$ws.Range("L174").Value = 0.136
$ws.Range("M174").Value = 0.30156
$ws.Range("N174").Value = 0.42132
$ws.Range("O174").Value = 1799980.004
$ws.Range("P174").Value = 5547896.885
$ws.Range("Q174").Value = "Rangitikei District"
$ws.Range("R174").Value = "Rangitīkei-Turakina"
$ws.Range("S174").Value = "Coastal Rangitikei"
$ws.Range("T174").Value = "Rang_4d"
$ws.Range("U174").Value = "mg/L"

# Row 175
$ws.Range("A175").Value = "Rangitikei at d/s Riverlands"
$ws.Range("B175").Value = "E coli (>260)"
$ws.Range("C175").Value = "B"
$ws.Range("D175").Value = "2019 - 2023"
$ws.Range("E175").Value = "Impact"
$ws.Range("F175").Value = 120
$ws.Range("G175").Value = 462.752870400332
$ws.Range("H175").Value = 8885.66648321926
$ws.Range("I175").Value = 1236
$ws.Range("J175").Value = 10.3448275862069
$ws.Range("K175").Value = 27.5862068965517
$ws.Range("L175").Value = 100
$ws.Range("M175").Value = 366.88
$ws.Range("N175").Value = 644.6799999999999
$ws.Range("O175").Value = 1799980.004
$ws.Range("P175").Value = 5547896.885
$ws.Range("Q175").Value = "Rangitikei District"
$ws.Range("R175").Value = "Rangitīkei-Turakina"
$ws.Range("S175").Value = "Coastal Rangitikei"
$ws.Range("T175").Value = "Rang_4d"
$ws.Range("U175").Value = "% exceedances over 260/100 mL"

# Row 176
$ws.Range("A176").Value = "Rangitikei at d/s Riverlands"
$ws.Range("B176").Value = "E coli (>540)"
$ws.Range("C176").Value = "C"
$ws.Range("D176").Value = "2019 - 2023"
$ws.Range("E176").Value = "Impact"
$ws.Range("F176").Value = 120
$ws.Range("G176").Value = 462.752870400332
$ws.Range("H176").Value = 8885.66648321926
$ws.Range("I176").Value = 1236
$ws.Range("J176").Value = 10.3448275862069
$ws.Range("K176").Value = 27.5862068965517
$ws.Range("L176").Value = 100
$ws.Range("M176").Value = 366.88
$ws.Range("N176").Value = 644.6799999999999
$ws.Range("O176").Value = 1799980.004
$ws.Range("P176").Value = 5547896.885
$ws.Range("Q176").Value = "Rangitikei District"
$ws.Range("R176").Value = "Rangitīkei-Turakina"
$ws.Range("S176").Value = "Coastal Rangitikei"
$ws.Range("T176").Value = "Rang_4d"
$ws.Range("U176").Value = "% exceedances over 540/100 mL"

# Row 177
$ws.Range("A177").Value = "Rangitikei at d/s Riverlands"
$ws.Range("B177").Value = "E coli (Median)"
$ws.Range("C177").Value = "A"
$ws.Range("D177").Value = "2019 - 2023"
$ws.Range("E177").Value = "Impact"
$ws.Range("F177").Value = 120
$ws.Range("G177").Value = 462.752870400332
$ws.Range("H177").Value = 8885.66648321926
$ws.Range("I177").Value = 1236
$ws.Range("J177").Value = 10.3448275862069
$ws.Range("K177").Value = 27.5862068965517
$ws.Range("L177").Value = 100
$ws.Range("M177").Value = 366.88
$ws.Range("N177").Value = 644.6799999999999
$ws.Range("O177").Value = 1799980.004
$ws.Range("P177").Value = 5547896.885
$ws.Range("Q177").Value = "Rangitikei District"
$ws.Range("R177").Value = "Rangitīkei-Turakina"
$ws.Range("S177").Value = "Coastal Rangitikei"
$ws.Range("T177").Value = "Rang_4d"
$ws.Range("U177").Value = "E. coli/100 mL"

# Row 178
$ws.Range("A178").Value = "Rangitikei at d/s Riverlands"
$ws.Range("B178").Value = "E coli (95th Percentile)"
$ws.Range("C178").Value = "E"
$ws.Range("D178").Value = "2019 - 2023"
$ws.Range("E178").Value = "Impact"
$ws.Range("F178").Value = 120
$ws.Range("G178").Value = 462.752870400332
$ws.Range("H178").Value = 8885.66648321926
$ws.Range("I178").Value = 1236
$ws.Range("J178").Value = 10.3448275862069
$ws.Range("K178").Value = 27.5862068965517
$ws.Range("L178").Value = 100
$ws.Range("M178").Value = 366.88
$ws.Range("N178").Value = 644.6799999999999
$ws.Range("O178").Value = 1799980.004
$ws.Range("P178").Value = 5547896.885
$ws.Range("Q178").Value = "Rangitikei District"
$ws.Range("R178").Value = "Rangitīkei-Turakina"
$ws.Range("S178").Value = "Coastal Rangitikei"
$ws.Range("T178").Value = "Rang_4d"
$ws.Range("U178").Value = "E. coli/100 mL"

# Row 179
$ws.Range("A179").Value = "Rangitikei at d/s Riverlands"
$ws.Range("B179").Value = "Ammoniacal-N (95th Percentile)"
$ws.Range("C179").Value = "C"
$ws.Range("D179").Value = "2019 - 2023"
$ws.Range("E179").Value = "Impact"
$ws.Range("F179").Value = 0.04347
$ws.Range("G179").Value = 0.195548142271213
$ws.Range("H179").Value = 2.9625657674792
$ws.Range("I179").Value = 0.8297600000000001
$ws.Range("J179").Value = ""
$ws.Range("K179").Value = ""
$ws.Range("L179").Value = 0.04418
$ws.Range("M179").Value = 0.29608
$ws.Range("N179").Value = 0.56168
$ws.Range("O179").Value = 1799980.004
$ws.Range("P179").Value = 5547896.885
$ws.Range("Q179").Value = "Rangitikei District"
$ws.Range("R179").Value = "Rangitīkei-Turakina"
$ws.Range("S179").Value = "Coastal Rangitikei"
$ws.Range("T179").Value = "Rang_4d"
$ws.Range("U179").Value = "mg NH4-N/L"

# Row 180
$ws.Range("A180").Value = "Rangitikei at d/s Riverlands"
$ws.Range("B180").Value = "Ammoniacal-N (Median)"
$ws.Range("C180").Value = "B"
$ws.Range("D180").Value = "2019 - 2023"
$ws.Range("E180").Value = "Impact"
$ws.Range("F180").Value = 0.04347
$ws.Range("G180").Value = 0.195548142271213
$ws.Range("H180").Value = 2.9625657674792
$ws.Range("I180").Value = 0.8297600000000001
$ws.Range("J180").Value = ""
$ws.Range("K180").Value = ""
$ws.Range("L180").Value = 0.04418
$ws.Range("M180").Value = 0.29608
$ws.Range("N180").Value = 0.56168
$ws.Range("O180").Value = 1799980.004
$ws.Range("P180").Value = 5547896.885
$ws.Range("Q180").Value = "Rangitikei District"
$ws.Range("R180").Value = "Rangitīkei-Turakina"
$ws.Range("S180").Value = "Coastal Rangitikei"
$ws.Range("T180").Value = "Rang_4d"
$ws.Range("U180").Value = "mg NH4-N/L"

# Row 181
$ws.Range("A181").Value = "Rangitikei at d/s Riverlands"
$ws.Range("B181").Value = "Nitrate-N (95th Percentile)"
$ws.Range("C181").Value = "A"
$ws.Range("D181").Value = "2019 - 2023"
$ws.Range("E181").Value = "Impact"
$ws.Range("F181").Value = 0.3355
$ws.Range("G181").Value = 0.430034482758621
$ws.Range("H181").Value = 2.68
$ws.Range("I181").Value = 1.1124
$ws.Range("J181").Value = ""
$ws.Range("K181").Value = ""
$ws.Range("L181").Value = 0.194
$ws.Range("M181").Value = 0.68412
$ws.Range("N181").Value = 0.9133599999999999
$ws.Range("O181").Value = 1799980.004
$ws.Range("P181").Value = 5547896.885
$ws.Range("Q181").Value = "Rangitikei District"
$ws.Range("R181").Value = "Rangitīkei-Turakina"
$ws.Range("S181").Value = "Coastal Rangitikei"
$ws.Range("T181").Value = "Rang_4d"
$ws.Range("U181").Value = "mg NO3-N/L"

# Row 182
$ws.Range("A182").Value = "Rangitikei at d/s Riverlands"
$ws.Range("B182").Value = "Nitrate-N (Median)"
$ws.Range("C182").Value = "A"
$ws.Range("D182").Value = "2019 - 2023"
$ws.Range("E182").Value = "Impact"
$ws.Range("F182").Value = 0.3355
$ws.Range("G182").Value = 0.430034482758621
$ws.Range("H182").Value = 2.68
$ws.Range("I182").Value = 1.1124
$ws.Range("J182").Value = ""
$ws.Range("K182").Value = ""
$ws.Range("L182").Value = 0.194
$ws.Range("M182").Value = 0.68412
$ws.Range("N182").Value = 0.9133599999999999
$ws.Range("O182").Value = 1799980.004
$ws.Range("P182").Value = 5547896.885
$ws.Range("Q182").Value = "Rangitikei District"
$ws.Range("R182").Value = "Rangitīkei-Turakina"
$ws.Range("S182").Value = "Coastal Rangitikei"
$ws.Range("T182").Value = "Rang_4d"
$ws.Range("U182").Value = "mg NO3-N/L"

# Row 183
$ws.Range("A183").Value = "Rangitikei at d/s Riverlands"
$ws.Range("B183").Value = "Soluble Inorganic Nitrogen (95th Percentile)"
$ws.Range("C183").Value = ""
$ws.Range("D183").Value = "2019 - 2023"
$ws.Range("E183").Value = "Impact"
$ws.Range("F183").Value = 0.4
$ws.Range("G183").Value = 0.569146551724138
$ws.Range("H183").Value = 3.008
$ws.Range("I183").Value = 1.8368
$ws.Range("J183").Value = ""
$ws.Range("K183").Value = ""
$ws.Range("L183").Value = 0.26
$ws.Range("M183").Value = 0.86556
$ws.Range("N183").Value = 1.4374
$ws.Range("O183").Value = 1799980.004
$ws.Range("P183").Value = 5547896.885
$ws.Range("Q183").Value = "Rangitikei District"
$ws.Range("R183").Value = "Rangitīkei-Turakina"
$ws.Range("S183").Value = "Coastal Rangitikei"
$ws.Range("T183").Value = "Rang_4d"
$ws.Range("U183").Value = "g/m3"

# Row 184
$ws.Range("A184").Value = "Rangitikei at d/s Riverlands"
$ws.Range("B184").Value = "Soluble Inorganic Nitrogen (Median)"
$ws.Range("C184").Value = ""
$ws.Range("D184").Value = "2019 - 2023"
$ws.Range("E184").Value = "Impact"
$ws.Range("F184").Value = 0.4
$ws.Range("G184").Value = 0.569146551724138
$ws.Range("H184").Value = 3.008
$ws.Range("I184").Value = 1.8368
$ws.Range("J184").Value = ""
$ws.Range("K184").Value = ""
$ws.Range("L184").Value = 0.26
$ws.Range("M184").Value = 0.86556
$ws.Range("N184").Value = 1.4374
$ws.Range("O184").Value = 1799980.004
$ws.Range("P184").Value = 5547896.885
$ws.Range("Q184").Value = "Rangitikei District"
$ws.Range("R184").Value = "Rangitīkei-Turakina"
$ws.Range("S184").Value = "Coastal Rangitikei"
$ws.Range("T184").Value = "Rang_4d"
$ws.Range("U184").Value = "g/m3"

# Row 185
$ws.Range("A185").Value = "Rangitikei at d/s Riverlands"
$ws.Range("B185").Value = "Total Nitrogen (95th Percentile)"
$ws.Range("C185").Value = ""
$ws.Range("D185").Value = "2019 - 2023"
$ws.Range("E185").Value = "Impact"
$ws.Range("F185").Value = 0.5600000000000001
$ws.Range("G185").Value = 0.816551724137931
$ws.Range("H185").Value = 3.08
$ws.Range("I185").Value = 2.268
$ws.Range("J185").Value = ""
$ws.Range("K185").Value = ""
$ws.Range("L185").Value = 0.8100000000000001
$ws.Range("M185").Value = 1.5284
$ws.Range("N185").Value = 1.9316
$ws.Range("O185").Value = 1799980.004
$ws.Range("P185").Value = 5547896.885
$ws.Range("Q185").Value = "Rangitikei District"
$ws.Range("R185").Value = "Rangitīkei-Turakina"
$ws.Range("S185").Value = "Coastal Rangitikei"
$ws.Range("T185").Value = "Rang_4d"
$ws.Range("U185").Value = "g/m3"

# Row 186
$ws.Range("A186").Value = "Rangitikei at d/s Riverlands"
$ws.Range("B186").Value = "Total Nitrogen (Median)"
$ws.Range("C186").Value = ""
$ws.Range("D186").Value = "2019 - 2023"
$ws.Range("E186").Value = "Impact"
$ws.Range("F186").Value = 0.5600000000000001
$ws.Range("G186").Value = 0.816551724137931
$ws.Range("H186").Value = 3.08
$ws.Range("I186").Value = 2.268
$ws.Range("J186").Value = ""
$ws.Range("K186").Value = ""
$ws.Range("L186").Value = 0.8100000000000001
$ws.Range("M186").Value = 1.5284
$ws.Range("N186").Value = 1.9316
$ws.Range("O186").Value = 1799980.004
$ws.Range("P186").Value = 5547896.885
$ws.Range("Q186").Value = "Rangitikei District"
$ws.Range("R186").Value = "Rangitīkei-Turakina"
$ws.Range("S186").Value = "Coastal Rangitikei"
$ws.Range("T186").Value = "Rang_4d"
$ws.Range("U186").Value = "g/m3"

# Row 187
$ws.Range("A187").Value = "Rangitikei at d/s Riverlands"
$ws.Range("B187").Value = "Total Phosphorus (95th Percentile)"
$ws.Range("C187").Value = ""
$ws.Range("D187").Value = "2019 - 2023"
$ws.Range("E187").Value = "Impact"
$ws.Range("F187").Value = 0.1445
$ws.Range("G187").Value = 0.245431034482759
$ws.Range("H187").Value = 1.41
$ws.Range("I187").Value = 0.8624000000000001
$ws.Range("J187").Value = ""
$ws.Range("K187").Value = ""
$ws.Range("L187").Value = 0.193
$ws.Range("M187").Value = 0.4698
$ws.Range("N187").Value = 0.75744
$ws.Range("O187").Value = 1799980.004
$ws.Range("P187").Value = 5547896.885
$ws.Range("Q187").Value = "Rangitikei District"
$ws.Range("R187").Value = "Rangitīkei-Turakina"
$ws.Range("S187").Value = "Coastal Rangitikei"
$ws.Range("T187").Value = "Rang_4d"
$ws.Range("U187").Value = "g/m3"

# Row 188
$ws.Range("A188").Value = "Rangitikei at d/s Riverlands"
$ws.Range("B188").Value = "Total Phosphorus (Median)"
$ws.Range("C188").Value = ""
$ws.Range("D188").Value = "2019 - 2023"
$ws.Range("E188").Value = "Impact"
$ws.Range("F188").Value = 0.1445
$ws.Range("G188").Value = 0.245431034482759
$ws.Range("H188").Value = 1.41
$ws.Range("I188").Value = 0.8624000000000001
$ws.Range("J188").Value = ""
$ws.Range("K188").Value = ""
$ws.Range("L188").Value = 0.193
$ws.Range("M188").Value = 0.4698
$ws.Range("N188").Value = 0.75744
$ws.Range("O188").Value = 1799980.004
$ws.Range("P188").Value = 5547896.885
$ws.Range("Q188").Value = "Rangitikei District"
$ws.Range("R188").Value = "Rangitīkei-Turakina"
$ws.Range("S188").Value = "Coastal Rangitikei"
$ws.Range("T188").Value = "Rang_4d"
$ws.Range("U188").Value = "g/m3"

